# Add a new JAX-WS demo entry ("webservice/jaxwsdemo1") after the existing
# "udpdemo1" section, and move the trailing "_GoBack" bookmark onto the new
# description paragraph.
#
# Before:
#   ...<w:p> ...udp协议发送中文与接收中文<bookmarkStart/><bookmarkEnd/></w:p>
#   <w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr></w:p>   <- trailing empty paragraph
#
# After:
#   ...<w:p> ...udp协议发送中文与接收中文</w:p>                               <- bookmark removed
#   <w:p/>                                                                  <- blank spacer paragraph
#   <w:p> (Heading 2, numbered) webservice/jaxwsdemo1 </w:p>
#   <w:p> 使用原生java代码编写一个最基础的webservice服务端并生成客户端调用<bookmarkStart/><bookmarkEnd/></w:p>

$d = $word.ActiveDocument

# Locate the "...udp协议发送中文与接收中文" paragraph (carries the _GoBack
# bookmark today) by its distinctive text rather than assuming a fixed
# index, then grab the paragraph right after it (the trailing blank
# paragraph that currently precedes <w:sectPr>).
$pPrev = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text -like "*协议发送中文与接收中文*") {
        $pPrev = $candidate
    }
}
$pLast = $pPrev.Next()

# Whole-range replace covering both the closing paragraph and the trailing
# blank paragraph so the old bookmark and old formatting are fully replaced
# by the freshly authored content in one shot.
$r = $d.Range($pPrev.Range.Start, $pLast.Range.End)

$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$xml = ""
$xml += "<w:p $w>"
$xml +=   "<w:r><w:rPr><w:rFonts w:hint=""eastAsia""/></w:rPr><w:t>使用</w:t></w:r>"
$xml +=   "<w:proofErr w:type=""spellStart""/>"
$xml +=   "<w:r><w:rPr><w:rFonts w:hint=""eastAsia""/></w:rPr><w:t>udp</w:t></w:r>"
$xml +=   "<w:proofErr w:type=""spellEnd""/>"
$xml +=   "<w:r><w:rPr><w:rFonts w:hint=""eastAsia""/></w:rPr><w:t>协议发送中文与接收中文</w:t></w:r>"
$xml += "</w:p>"

$xml += "<w:p $w/>"

$xml += "<w:p $w>"
$xml +=   "<w:pPr>"
$xml +=     "<w:pStyle w:val=""2""/>"
$xml +=     "<w:numPr><w:ilvl w:val=""0""/><w:numId w:val=""1""/></w:numPr>"
$xml +=     "<w:rPr><w:rFonts w:hint=""eastAsia""/></w:rPr>"
$xml +=   "</w:pPr>"
$xml +=   "<w:r><w:rPr><w:rFonts w:hint=""eastAsia""/></w:rPr><w:lastRenderedPageBreak/><w:t>w</w:t></w:r>"
$xml +=   "<w:r><w:t>ebservice</w:t></w:r>"
$xml +=   "<w:r><w:rPr><w:rFonts w:hint=""eastAsia""/></w:rPr><w:t>/</w:t></w:r>"
$xml +=   "<w:r><w:t>jaxwsdemo1</w:t></w:r>"
$xml += "</w:p>"

$xml += "<w:p $w>"
$xml +=   "<w:r><w:rPr><w:rFonts w:hint=""eastAsia""/></w:rPr><w:t>使用原生java代码编写一个最基础的webservice服务端并生成客户端调用</w:t></w:r>"
$xml +=   "<w:bookmarkStart w:id=""0"" w:name=""_GoBack""/>"
$xml +=   "<w:bookmarkEnd w:id=""0""/>"
$xml += "</w:p>"

$null = $r.InsertXML($xml)
